{"js": "// Add a new \"Rule:\" bullet (top level) followed by a sub-bullet describing\n// the rule, mirroring the author's addition after the \"Write process\n// specifications for Process 1,3,4,6\" item.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// New top-level list item: \"Rule:\"\nconst ruleHeading = lastParagraph.insertParagraph(\"Rule:\", Word.InsertLocation.after);\nruleHeading.listItem.level = 0;\n\n// New second-level list item with the rule text.\nconst ruleBody = ruleHeading.insertParagraph(\"Any documents must be uploaded before 11h45pm.\", Word.InsertLocation.after);\nruleBody.listItem.level = 1;\n\nawait context.sync();\n", "ps1": "# Add a new \"Rule:\" bullet (top level) followed by a sub-bullet describing\n# the rule, mirroring the author's addition after the \"Write process\n# specifications for Process 1,3,4,6\" item.\n\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n\n# New top-level list item: \"Rule:\"\n$lastParagraph.Range.InsertParagraphAfter()\n$d = $word.ActiveDocument\n$ruleHeading = $d.Paragraphs.Last\n$ruleHeading.Range.Text = \"Rule:\"\n$ruleHeading.Range.ListFormat.ListLevelNumber = 1\n\n# New second-level list item with the rule text.\n$ruleHeading.Range.InsertParagraphAfter()\n$d = $word.ActiveDocument\n$ruleBody = $d.Paragraphs.Last\n$ruleBody.Range.Text = \"Any documents must be uploaded before 11h45pm.\"\n$ruleBody.Range.ListFormat.ListLevelNumber = 2\n"}
